$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 177, shifting the existing rows (177-184) down to (178-185)
$ws.Rows.Item(177).Insert()

# Populate the new row 177 with the new weekly record
$ws.Cells.Item(177, 1).Value = 11
$ws.Cells.Item(177, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(177, 3).Value = "Bíobío"
$ws.Cells.Item(177, 4).Value = 44747
$ws.Cells.Item(177, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(177, 5).Value = 8
$ws.Cells.Item(177, 6).Value = 100112003
$ws.Cells.Item(177, 7).Value = "Ajo"
$ws.Cells.Item(177, 8).Value = "Chino"
$ws.Cells.Item(177, 9).Value = "Primera"
$ws.Cells.Item(177, 10).Value = 350
$ws.Cells.Item(177, 11).Value = 18000
$ws.Cells.Item(177, 12).Value = 19000
$ws.Cells.Item(177, 13).Value = 18571
$ws.Cells.Item(177, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(177, 15).Value = "China"
$ws.Cells.Item(177, 16).Value = 1857
$ws.Cells.Item(177, 17).Value = 10
$ws.Cells.Item(177, 18).Value = "Hortaliza"
